$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "23.424.55"
$ws.Range("E2").Value = "  -1.36%  "

# Row 3
$ws.Range("D3").Value = "1.643.26"
$ws.Range("E3").Value = "  -0.72%  "

# Row 4
$origStyle = $ws.Range("D4").Style
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("D4").Style = $origStyle
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$origStyle = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "1.001"
$ws.Range("D5").Style = $origStyle
$ws.Range("E5").Value = "  +0.08%  "

# Row 6
$origStyle = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "298.72"
$ws.Range("D6").Style = $origStyle
$ws.Range("E6").Value = "  -1.73%  "

# Row 7
$origStyle = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3783"
$ws.Range("D7").Style = $origStyle
$ws.Range("E7").Value = "  -1.06%  "

# Row 8
$ws.Range("B8").Value = "OKB"
$ws.Range("C8").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$origStyle = $ws.Range("D8").Style
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "50.06"
$ws.Range("D8").Style = $origStyle
$ws.Range("E8").Value = "  -2.22%  "

# Row 9
$ws.Range("B9").Value = "Cardano"
$ws.Range("C9").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$origStyle = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.3521"
$ws.Range("D9").Style = $origStyle
$ws.Range("E9").Value = "  -2.53%  "

# Row 10
$origStyle = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.08068"
$ws.Range("D10").Style = $origStyle
$ws.Range("E10").Value = "  -1.81%  "

# Row 11
$origStyle = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.209"
$ws.Range("D11").Style = $origStyle
$ws.Range("E11").Value = "  -3.66%  "

# Row 12
$origStyle = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"
$ws.Range("D12").Style = $origStyle
$ws.Range("E12").Value = "  +0.07%  "

# Row 13
$origStyle = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.01"
$ws.Range("D13").Style = $origStyle
$ws.Range("E13").Value = "  -3.15%  "

# Row 14
$origStyle = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.373"
$ws.Range("D14").Style = $origStyle
$ws.Range("E14").Value = "  -2.70%  "

# Row 15
$origStyle = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.293"
$ws.Range("D15").Style = $origStyle
$ws.Range("E15").Value = "  -2.13%  "

# Row 16
$origStyle = $ws.Range("D16").Style
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.00001199"
$ws.Range("D16").Style = $origStyle
$ws.Range("E16").Value = "  -3.36%  "

# Row 17
$ws.Range("D17").Value = "1.635.81"
$ws.Range("E17").Value = "  +0.01%  "

# Row 18
$origStyle = $ws.Range("D18").Style
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "96.59"
$ws.Range("D18").Style = $origStyle
$ws.Range("E18").Value = "  -1.26%  "

# Row 19
$origStyle = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06970"
$ws.Range("D19").Style = $origStyle
$ws.Range("E19").Value = "  -0.19%  "

# Row 20
$origStyle = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.742"
$ws.Range("D20").Style = $origStyle
$ws.Range("E20").Value = "  -0.73%  "

# Row 21
$origStyle = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.33"
$ws.Range("D21").Style = $origStyle
$ws.Range("E21").Value = "  -2.51%  "

# Row 22
$ws.Range("E22").Value = "  +0.00%  "

# Row 23
$origStyle = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "12.37"
$ws.Range("D23").Style = $origStyle
$ws.Range("E23").Value = "  -2.97%  "

# Row 24
$ws.Range("D24").Value = "23.448.71"
$ws.Range("E24").Value = "  -1.28%  "

# Row 25
$origStyle = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.491"
$ws.Range("D25").Style = $origStyle
$ws.Range("E25").Value = "  -2.89%  "

# Row 26
$origStyle = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.887"
$ws.Range("D26").Style = $origStyle
$ws.Range("E26").Value = "  -6.51%  "

# Row 27
$origStyle = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.83"
$ws.Range("D27").Style = $origStyle
$ws.Range("E27").Value = "  -2.32%  "

# Row 28
$origStyle = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "152.84"
$ws.Range("D28").Style = $origStyle
$ws.Range("E28").Value = "  +1.29%  "

# Row 29
$origStyle = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.208"
$ws.Range("D29").Style = $origStyle
$ws.Range("E29").Value = "  -0.29%  "

# Row 30
$origStyle = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "132.21"
$ws.Range("D30").Style = $origStyle
$ws.Range("E30").Value = "  -2.04%  "

# Row 31
$ws.Range("D31").Value = "1.825.49"
$ws.Range("E31").Value = "  +0.27%  "

# Row 32
$origStyle = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.883"
$ws.Range("D32").Style = $origStyle
$ws.Range("E32").Value = "  -0.57%  "

# Row 33
$origStyle = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.138"
$ws.Range("D33").Style = $origStyle
$ws.Range("E33").Value = "  -1.42%  "

# Row 34
$origStyle = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "11.48"
$ws.Range("D34").Style = $origStyle
$ws.Range("E34").Value = "  -3.45%  "

# Row 35
$origStyle = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9848"
$ws.Range("D35").Style = $origStyle
$ws.Range("E35").Value = "  -9.06%  "

# Row 36
$origStyle = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02692"
$ws.Range("D36").Style = $origStyle
$ws.Range("E36").Value = "  -5.03%  "

# Row 37
$origStyle = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.08728"
$ws.Range("D37").Style = $origStyle
$ws.Range("E37").Value = "  -1.16%  "

# Row 38
$origStyle = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2434"
$ws.Range("D38").Style = $origStyle
$ws.Range("E38").Value = "  -3.43%  "

# Row 39
$origStyle = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.917"
$ws.Range("D39").Style = $origStyle
$ws.Range("E39").Value = "  -3.97%  "

# Row 40
$origStyle = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.06793"
$ws.Range("D40").Style = $origStyle
$ws.Range("E40").Value = "  -5.47%  "

# Row 41
$origStyle = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.85"
$ws.Range("D41").Style = $origStyle
$ws.Range("E41").Value = "  -1.81%  "

# Row 42
$origStyle = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.6848"
$ws.Range("D42").Style = $origStyle
$ws.Range("E42").Value = "  -3.20%  "

# Row 43
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$origStyle = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "15.67"
$ws.Range("D43").Style = $origStyle
$ws.Range("E43").Value = "  -1.78%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$origStyle = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.288"
$ws.Range("D44").Style = $origStyle
$ws.Range("E44").Value = "  -3.94%  "

# Row 45
$ws.Range("E45").Value = "  +0.11%  "

# Row 46
$origStyle = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.6324"
$ws.Range("D46").Style = $origStyle
$ws.Range("E46").Value = "  -3.46%  "

# Row 47
$origStyle = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.246"
$ws.Range("D47").Style = $origStyle
$ws.Range("E47").Value = "  -3.67%  "

# Row 48
$origStyle = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.899"
$ws.Range("D48").Style = $origStyle
$ws.Range("E48").Value = "  -1.55%  "

# Row 49
$origStyle = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.07717"
$ws.Range("D49").Style = $origStyle
$ws.Range("E49").Value = "  -3.27%  "

# Row 50
$origStyle = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "126.84"
$ws.Range("D50").Style = $origStyle
$ws.Range("E50").Value = "  -1.40%  "

# Row 51
$origStyle = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.143"
$ws.Range("D51").Style = $origStyle
$ws.Range("E51").Value = "  -4.33%  "
